{"js": "// Find the \"Address: ...\" paragraph and insert a new paragraph right after\n// it containing \"Website: hholben.github.io\". The new paragraph should\n// inherit the same paragraph/run formatting as the Address paragraph\n// (Times New Roman, 12pt, centered, black), matching how Word behaves when\n// you press Enter at the end of that line and type the new text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet addressPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(\"Address:\") !== -1) {\n    addressPara = p;\n    break;\n  }\n}\n\nif (!addressPara) {\n  throw new Error(\"Could not find the 'Address:' paragraph.\");\n}\n\n// Inserting a paragraph *after* the Address paragraph duplicates its\n// paragraph mark formatting (pPr/rPr), exactly like pressing Enter at the\n// end of that paragraph in Word.\nconst newPara = addressPara.insertParagraph(\"Website: hholben.github.io\", \"After\");\n\nawait context.sync();\n", "ps1": "# Find the \"Address: ...\" paragraph and insert a brand-new paragraph right\n# after it containing \"Website: hholben.github.io\". Word's\n# InsertParagraphAfter() duplicates the paragraph mark formatting (pPr/rPr)\n# of the Address paragraph onto the new paragraph, just like pressing Enter\n# at the end of that line and typing the new text.\n\n$d = $word.ActiveDocument\n\n$findRange = $d.Content\n$null = $findRange.Find.Execute(\"Address:\")\nif (-not $findRange.Find.Found) {\n    throw \"Could not find the 'Address:' paragraph.\"\n}\n\n$addressPara = $findRange.Paragraphs(1)\n$addressRange = $addressPara.Range\n$null = $addressRange.InsertParagraphAfter()\n\n$newPara = $addressPara.Next()\n$newPara.Range.Text = \"Website: hholben.github.io\"\n"}
